$wb = $excel.ActiveWorkbook

# --- Sheet "Logs": fix row 34 (remove stray "nan" answer) and append row 35 ---
$wsLogs = $wb.Worksheets.Item("Logs")

# Row 34: E34 contained the literal text "nan" - clear it so the cell is removed.
$wsLogs.Range("E34").Value = ""

# New row 35: an incoming "Inlogproblemen" ticket that has now been answered.
$wsLogs.Range("A35").Value = "Inlogproblemen"
$wsLogs.Range("B35").Value = "mailmind.test@zohomail.eu"
$wsLogs.Range("C35").Value = "Ik kan niet meer inloggen op mijn account. Kunnen jullie helpen?"
$wsLogs.Range("D35").Value = "IT / Technisch probleem"
$wsLogs.Range("E35").Value = "Beste klant,
Bedankt voor je bericht. Om je verder te kunnen helpen met het inlogprobleem, hebben we wat meer informatie nodig. Zou je alsjeblieft je gebruikersnaam en eventuele foutmeldingen die je hebt ontvangen kunnen delen? Op die manier kunnen we het probleem zo snel mogelijk voor je oplossen.
Met vriendelijke groet,
[Bedrijfsnaam] E-mailassistent"
$wsLogs.Range("F35").Value = "2025-06-24 21:40:16"
$wsLogs.Range("G35").Value = "Ja"

# The conditional formatting ranges need to grow from row 34 to row 35.
$fcCategorie = $wsLogs.Range("D2:D34").FormatConditions
$fcCategorie.Item(1).ModifyAppliesToRange($wsLogs.Range("D2:D35"))

$fcBeantwoord = $wsLogs.Range("G2:G34").FormatConditions
$fcBeantwoord.Item(1).ModifyAppliesToRange($wsLogs.Range("G2:G35"))

# --- Sheet "Dashboard": the counts for "Bestelling / Levering" and
#     "IT / Technisch probleem" swapped order/rank because of the new ticket ---
$wsDash = $wb.Worksheets.Item("Dashboard")
$wsDash.Range("A4").Value = "IT / Technisch probleem"
$wsDash.Range("A5").Value = "Bestelling / Levering"
$wsDash.Range("B5").Value = 4
